# Update the Fgf6-Fgfr3 NATMI output sheet with refreshed TPM-derived values
# and drop the rows for clusters that no longer have data (Neutrophils,
# Resolving-Mac) as well as the old Inflammatory-Mac row, leaving only the
# ECs, FAPs and (new) MuSCs target-cluster rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Target cluster: ECs) -----------------------------------------
$ws.Range("M2").Value = 7.6704545
$ws.Range("N2").Value = 15.340909
$ws.Range("O2").Value = 0.8278663930876066
$ws.Range("P2").Value = 0.7913005936208135
$ws.Range("Q2").Value = 0.03395965888966667
$ws.Range("R2").Value = 0.203757953338
$ws.Range("S2").Value = 0.8278663930876066
$ws.Range("T2").Value = 0.7913005936208135

# --- Row 3 (Target cluster: FAPs) -----------------------------------------
$ws.Range("O3").Value = 0.0924194920655273
$ws.Range("P3").Value = 0.1325061620042962
$ws.Range("S3").Value = 0.0924194920655273
$ws.Range("T3").Value = 0.1325061620042962

# --- Row 4 (Target cluster: was Inflammatory-Mac, now MuSCs) --------------
$ws.Range("D4").Value = "MuSCs"
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7385775
$ws.Range("N4").Value = 1.477155
$ws.Range("O4").Value = 0.0797141148468662
$ws.Range("P4").Value = 0.07619324437489022
$ws.Range("Q4").Value = 0.003269928785
$ws.Range("R4").Value = 0.01961957271
$ws.Range("S4").Value = 0.0797141148468662
$ws.Range("T4").Value = 0.07619324437489022

# --- Remove the now-obsolete rows (MuSCs, Neutrophils, Resolving-Mac) -----
$ws.Range("A5:T7").Delete() | Out-Null
